$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = 155000
$ws.Range("B25").Value = 66
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 'Santa Justa - Miraflores - Cruz Roja'
$ws.Range("F25").Value = 'Arroyo - Santa Justa'

$ws.Range("A26").Value = 205000
$ws.Range("B26").Value = 75
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 'La Palmera - Los Bermejales'
$ws.Range("F26").Value = 'Bami - Pineda'

$ws.Range("A27").Value = 280000
$ws.Range("B27").Value = 83
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 'Centro'
$ws.Range("F27").Value = 'Arenal - Museo - Tetuán'

$ws.Range("A28").Value = 230000
$ws.Range("B28").Value = 68
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 'Triana'
$ws.Range("F28").Value = 'López de Gomara'

$ws.Range("A29").Value = 320000
$ws.Range("B29").Value = 108
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = 'Triana'
$ws.Range("F29").Value = 'Ronda de Triana-Patrocinio-Turruñuelo'

$ws.Range("A30").Value = 280000
$ws.Range("B30").Value = 80
$ws.Range("C30").Value = 3
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = 'Nervión'
$ws.Range("F30").Value = 'Gran Plaza - Marqués de Pickman - Ciudad Jardín'

$ws.Range("A31").Value = 359000
$ws.Range("B31").Value = 102
$ws.Range("C31").Value = 3
$ws.Range("D31").Value = 2
$ws.Range("E31").Value = 'Nervión'
$ws.Range("F31").Value = 'Nervión'

$ws.Range("A32").Value = 339000
$ws.Range("B32").Value = 95
$ws.Range("C32").Value = 3
$ws.Range("D32").Value = 2
$ws.Range("E32").Value = 'Triana'
$ws.Range("F32").Value = 'Calle Betis - Pagés del Corro'

$ws.Range("A33").Value = 550000
$ws.Range("B33").Value = 124
$ws.Range("C33").Value = 3
$ws.Range("D33").Value = 2
$ws.Range("E33").Value = 'Centro'
$ws.Range("F33").Value = 'Arenal - Museo - Tetuán'

$ws.Range("A34").Value = 789000
$ws.Range("B34").Value = 286
$ws.Range("C34").Value = 6
$ws.Range("D34").Value = 5
$ws.Range("E34").Value = 'Nervión'
$ws.Range("F34").Value = 'Buhaira - Huerta del Rey'

$ws.Range("A35").Value = 170000
$ws.Range("B35").Value = 66
$ws.Range("C35").Value = 2
$ws.Range("D35").Value = 1
$ws.Range("E35").Value = 'Santa Justa - Miraflores - Cruz Roja'
$ws.Range("F35").Value = 'Arroyo - Santa Justa'

$ws.Range("A36").Value = 330000
$ws.Range("B36").Value = 106
$ws.Range("C36").Value = 3
$ws.Range("D36").Value = 2
$ws.Range("E36").Value = 'Los Remedios'
$ws.Range("F36").Value = 'Ramón de Carranza - Madre Rafols'

$ws.Range("A37").Value = 250000
$ws.Range("B37").Value = 73
$ws.Range("C37").Value = 3
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = 'Centro'
$ws.Range("F37").Value = 'Puerta Carmona-Puerta Osario-Amador de los Ríos'

$ws.Range("A38").Value = 300000
$ws.Range("B38").Value = 100
$ws.Range("C38").Value = 2
$ws.Range("D38").Value = 1
$ws.Range("E38").Value = 'Centro'
$ws.Range("F38").Value = 'San Vicente'

$ws.Range("A39").Value = 289000
$ws.Range("B39").Value = 82
$ws.Range("C39").Value = 3
$ws.Range("D39").Value = 2
$ws.Range("E39").Value = 'Triana'
$ws.Range("F39").Value = 'Ronda de Triana-Patrocinio-Turruñuelo'

$ws.Range("A40").Value = 499000
$ws.Range("B40").Value = 189
$ws.Range("C40").Value = 5
$ws.Range("D40").Value = 2
$ws.Range("E40").Value = 'Los Remedios'
$ws.Range("F40").Value = 'Ramón de Carranza - Madre Rafols'

$ws.Range("A41").Value = 349000
$ws.Range("B41").Value = 126
$ws.Range("C41").Value = 3
$ws.Range("D41").Value = 2
$ws.Range("E41").Value = 'Triana'
$ws.Range("F41").Value = 'Ronda de Triana-Patrocinio-Turruñuelo'

$ws.Range("A42").Value = 295000
$ws.Range("B42").Value = 158
$ws.Range("C42").Value = 3
$ws.Range("D42").Value = 1
$ws.Range("E42").Value = 'Triana'
$ws.Range("F42").Value = 'López de Gomara'

$ws.Range("A43").Value = 380000
$ws.Range("B43").Value = 189
$ws.Range("C43").Value = 6
$ws.Range("D43").Value = 3
$ws.Range("E43").Value = 'Nervión'
$ws.Range("F43").Value = 'Nervión'

$ws.Range("A44").Value = 330000
$ws.Range("B44").Value = 140
$ws.Range("C44").Value = 4
$ws.Range("D44").Value = 2
$ws.Range("E44").Value = 'Nervión'
$ws.Range("F44").Value = 'Nervión'

$ws.Range("A45").Value = 372000
$ws.Range("B45").Value = 99
$ws.Range("C45").Value = 2
$ws.Range("D45").Value = 1
$ws.Range("E45").Value = 'Centro'
$ws.Range("F45").Value = 'San Vicente'

$ws.Range("A46").Value = 368000
$ws.Range("B46").Value = 90
$ws.Range("C46").Value = 2
$ws.Range("D46").Value = 2
$ws.Range("E46").Value = 'Centro'
$ws.Range("F46").Value = 'San Vicente'

$ws.Range("A47").Value = 290000
$ws.Range("B47").Value = 114
$ws.Range("C47").Value = 3
$ws.Range("D47").Value = 2
$ws.Range("E47").Value = 'Nervión'
$ws.Range("F47").Value = 'Luis Montoto - Santa Justa'

$ws.Range("A48").Value = 1495000
$ws.Range("B48").Value = 346
$ws.Range("C48").Value = 5
$ws.Range("D48").Value = 3
$ws.Range("E48").Value = 'Centro'
$ws.Range("F48").Value = 'Santa Cruz - Alfalfa'

$ws.Range("A49").Value = 132260
$ws.Range("B49").Value = 114
$ws.Range("C49").Value = 4
$ws.Range("D49").Value = 2
$ws.Range("E49").Value = 'Cerro Amate'
$ws.Range("F49").Value = 'Santa Aurelia'

$ws.Range("A50").Value = 450000
$ws.Range("B50").Value = 114
$ws.Range("C50").Value = 3
$ws.Range("D50").Value = 2
$ws.Range("E50").Value = 'Nervión'
$ws.Range("F50").Value = 'Nervión'

$ws.Range("A51").Value = 575000
$ws.Range("B51").Value = 140
$ws.Range("C51").Value = 4
$ws.Range("D51").Value = 3
$ws.Range("E51").Value = 'Prado de San Sebastián - Felipe II - Bueno Monreal'
$ws.Range("F51").Value = 'Felipe II - Bueno Monreal'
